$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "D0.9125181743180528"
$ws.Range("A3").Value = "D0.5477877949969697"
$ws.Range("A4").Value = "D0.6037163037954872"

$ws.Range("D2").Value = "Fri, 23 Dec 2022 23:19:07 -0800"
$ws.Range("D3").Value = "Fri, 23 Dec 2022 23:19:07 -0800"
$ws.Range("D4").Value = "Fri, 23 Dec 2022 23:19:07 -0800"
